# Insert a new slide at position 22: "构建者的几次角色转换"
# This pushes the former slide 22 ("前置知识和技能" ... ) and everything after it
# down by one position, which matches the target presentation's new sldIdLst.

$p = $ppt.ActivePresentation

# Layout 13 of the (shared) slide master is "1_7*#标题和内容（一行标题）"
# -- the same Title+Body(idx=10) layout used by the neighbouring
# "为什么选中 LFS7.7-systemd" slide, and it is the layout that matches the
# two placeholders (title + body idx=10) added by this change.
$lay = $p.SlideMaster.CustomLayouts.Item(13)

$newSlide = $p.Slides.AddSlide(22, $lay)

# --- Title placeholder ---
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "构建者的几次角色转换"

# --- Body placeholder (idx=10) ---
$body = $newSlide.Shapes.Item(2)
$body.Left = 455612 / 12700.0
$body.Top = 933450 / 12700.0
$body.Width = 11293475 / 12700.0
$body.Height = 5267325 / 12700.0

$paragraphs = @(
    "最开始，宿主系统的 root 用户准备好宿主系统的软件环境，为宿主系统的 lfs 用户铺平了道路；",
    "然后，宿主系统的 lfs 用户构建好临时工具链，为宿主系统的 chrooted-root 用户铺平了道路；",
    "接着，chrooted-root 用户构建好基本系统软件、进行系统配置、编译并安装内核，为 LFS 目标系统的 root 用户铺平了道路；",
    "最后，宿主系统的 root 用户通过配置宿主系统的 GRUB 菜单为 LFS 目标系统在宿主系统启动阶段提供了一个入口；",
    "这样，最终 LFS 目标系统的 root 用户得以进入自己的系统！"
)

$body.TextFrame.TextRange.Text = [string]::Join("`r", $paragraphs)
